$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "U"
$ws.Range("D3").Value = "U"
$ws.Range("D4").Value = "U"
$ws.Range("D5").Value = "U"

$ws.Range("G12").Select()
